# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" / "Valor Mora" table (rows 16-22, cols E & F) gets its
# row order reversed: the period that used to be last (2107, with the
# different Valor Mora of 29260) now appears first, and the rest follow in
# descending order down to 2101 (which now carries the common 35112 value
# that the other rows had).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current Periodo Mora (col E) / Valor Mora (col F) pairs for
# rows 16-22 before touching anything.
$periods = @()
$values  = @()
for ($r = 16; $r -le 22; $r++) {
    $periods += $ws.Cells.Item($r, 5).Value2
    $values  += $ws.Cells.Item($r, 6).Value2
}

# Write them back in reverse order, so row 16 gets what used to be in row 22,
# row 17 gets what used to be in row 21, etc.
$n = $periods.Length
for ($i = 0; $i -lt $n; $i++) {
    $r = 16 + $i
    $srcIndex = $n - 1 - $i
    $ws.Cells.Item($r, 5).Value = $periods[$srcIndex]
    $ws.Cells.Item($r, 6).Value = $values[$srcIndex]
}
